$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Evaporator Temperature
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 30
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0.1

# Row 3 - Condenser Temperature
$ws.Range("D3").Value = 65
$ws.Range("E3").Value = 35
$ws.Range("H3").Value = "Discrete"

# Row 4 - Adiabatic Efficiency
$ws.Range("H4").Value = "Discrete"

# Row 5 - Capacity
$ws.Range("D5").Value = 20000
$ws.Range("E5").Value = 5000
$ws.Range("F5").Value = 100
$ws.Range("H5").Value = "Discrete"

# Row 2's H2 cell previously referenced the shared string "Continuous".
# Since "Continuous" is removed from the shared strings table entirely,
# H2 now resolves to "Discrete" (the string that took its place).
$ws.Range("H2").Value = "Discrete"
